$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the legend/caption text in A1 (rich text: bold title + regular body) ---

# Prepend "Locality " to the bold title run
$ws.Range("A1").Characters(1, 0).Text = "Locality "

# Insert " using locality" right before the closing period of the body paragraph
$full = $ws.Range("A1").Text
$idx = $full.IndexOf("elements.")
$insertPos = $idx + "elements".Length
$ws.Range("A1").Characters($insertPos + 1, 0).Text = " using locality"

# Re-apply run formatting so the title stays bold and the body stays regular
$full2 = $ws.Range("A1").Text
$titleLen = "Locality HPO genes discovered with networks built from accessions subsets".Length
$titleRun = $ws.Range("A1").Characters(1, $titleLen)
$titleRun.Font.Bold = $true
$titleRun.Font.Name = "Calibri"
$titleRun.Font.Size = 11
$bodyRun = $ws.Range("A1").Characters($titleLen + 1, $full2.Length - $titleLen)
$bodyRun.Font.Bold = $false
$bodyRun.Font.Name = "Calibri"
$bodyRun.Font.Size = 11

# The caption text grew by a couple words, so the wrapped title row needs to be taller
$ws.Rows.Item(1).RowHeight = 93.75

# --- Add a new "Total Ionome" summary row (row 21) ---

# Clone the row-4..19 cell style onto the new row label cell
[void]$ws.Range("A19").Copy()
[void]$ws.Range("A21").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A21").Value = "Total Ionome"
$ws.Range("B21").Formula = "=SUM(B4:B20)"
$ws.Range("C21:G21").Formula = "=SUM(C4:C20)"

# --- Move the active selection, matching the saved view state ---
[void]$ws.Range("E31").Select()
